# Update the '报名人数' (F column) values across all sheets to the latest
# scraped counts, per the gh-pages data refresh at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 843
$ws.Range("F3").Value = 13704
$ws.Range("F4").Value = 13485
$ws.Range("F5").Value = 1046
$ws.Range("F10").Value = 21
$ws.Range("F11").Value = 45
$ws.Range("F12").Value = 749
$ws.Range("F13").Value = 2135
$ws.Range("F14").Value = 77
$ws.Range("F15").Value = 85
$ws.Range("F16").Value = 69
$ws.Range("F17").Value = 111
$ws.Range("F19").Value = 509
$ws.Range("F21").Value = 380
$ws.Range("F22").Value = 311
$ws.Range("F23").Value = 68
$ws.Range("F24").Value = 820
$ws.Range("F25").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 47
$ws.Range("F6").Value = 157
$ws.Range("F7").Value = 1439
$ws.Range("F10").Value = 54

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 96

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 843
$ws.Range("F4").Value = 13704
$ws.Range("F5").Value = 13485
$ws.Range("F6").Value = 1046
$ws.Range("F11").Value = 21
$ws.Range("F12").Value = 45
$ws.Range("F13").Value = 749
$ws.Range("F14").Value = 47
$ws.Range("F16").Value = 2135
$ws.Range("F17").Value = 77
$ws.Range("F18").Value = 85
$ws.Range("F19").Value = 69
$ws.Range("F20").Value = 111
$ws.Range("F24").Value = 96
$ws.Range("F25").Value = 96
$ws.Range("F26").Value = 509
$ws.Range("F28").Value = 380
$ws.Range("F29").Value = 311
$ws.Range("F30").Value = 68
$ws.Range("F31").Value = 820
$ws.Range("F32").Value = 157
$ws.Range("F33").Value = 1439
$ws.Range("F36").Value = 72
$ws.Range("F37").Value = 54
